$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A4").Value = "User or Password is not valid"
$ws.Range("A4").HorizontalAlignment = -4108
$ws.Range("A4").VerticalAlignment = -4108

$ws.Range("A4").Select()
